$wb = $excel.ActiveWorkbook

# --- "Commands" sheet: insert a new row for the "playbackdropmode" command ---
$ws = $wb.Worksheets.Item("Commands")

# Insert a new blank row above the current row 109 ("playback(n,<bool>)"),
# pushing that row (and everything below it) down by one.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new command description.
$ws.Cells.Item(109, 2).Value = "playbackdropmode(<int>)"
$ws.Cells.Item(109, 3).Value = "sets playback DROP mode to 0: off, 1: time, 2: BT, 3: ET"

# Touch column D on the new row so it becomes part of the used range
# (mirrors the empty, styled <c r="D109"/> cell added in the target sheet).
$ws.Cells.Item(109, 4).NumberFormat = "General"

# --- view/selection bookkeeping (cosmetic, mirrors the saved cursor state) ---
$ws.Activate()
$ws.Range("A100").Select()
$ws.Range("A109").Select()

$ws1 = $wb.Worksheets.Item("Sliders")
$ws1.Range("B7").Select()

$ws.Activate()
